$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the EDIPI (C) values from rows 9-16 - clear the existing numeric EDIPI
$ws.Range("C9:C16").ClearContents()

# Add two new rows of data: row 17 and row 18
$ws.Range("A17").Value = "Fletcher"
$ws.Range("B17").Value = "Nicholas"
$ws.Range("C17").Value = 4453245321
$ws.Range("D17").Value = "ACTIVE DUTY"

$ws.Range("A18").Value = "Cena"
$ws.Range("B18").Value = "J"
$ws.Range("C18").Value = 5555555555
$ws.Range("D18").Value = "ACTIVE DUTY"

$ws.Range("E17").Value = "RandomCourse"
$ws.Range("E18").Value = "RandomCourse"

# Dates as raw serials (avoids auto-format detection creating a new
# numFmt); then copy the existing date-formatted cell's format over so
# the style index matches the rest of the F/G columns (s="1").
$ws.Range("F17").Value2 = 36526
$ws.Range("G17").Value2 = 36527
$ws.Range("F18").Value2 = 36526
$ws.Range("G18").Value2 = 36161

$ws.Range("F2").Copy()
$ws.Range("F17:F18").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("G17:G18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C9").Select()
